$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 151.375
$ws.Range("I11").Value = 151.375
$ws.Range("K11").Value = 151.375
$ws.Range("M11").Value = -11.375

$ws.Range("H98").Value = 1804.24
$ws.Range("I98").Value = 1317.75
$ws.Range("J98").Value = 3750.2
$ws.Range("K98").Value = 1317.75
$ws.Range("L98").Value = 3750.2
$ws.Range("M98").Value = 180.25
$ws.Range("N98").Value = -6746.2

$ws.Range("H116").Value = 7426.7896
$ws.Range("I116").Value = 12834.444
$ws.Range("J116").Value = 2559.9
$ws.Range("K116").Value = 12834.444
$ws.Range("L116").Value = 2559.9
$ws.Range("M116").Value = -9392.444
$ws.Range("N116").Value = -9443.9

$ws.Range("H122").Value = 1804.24
$ws.Range("I122").Value = 1317.75
$ws.Range("J122").Value = 3750.2
$ws.Range("K122").Value = 3953.25
$ws.Range("L122").Value = 11250.6
$ws.Range("M122").Value = -1503.25
$ws.Range("N122").Value = -16150.6

$ws.Range("H132").Value = 1116.6111
$ws.Range("I132").Value = 946.82355
$ws.Range("J132").Value = 4003
$ws.Range("K132").Value = 2840.47065
$ws.Range("L132").Value = 12009
$ws.Range("M132").Value = -310.4706499999998
$ws.Range("N132").Value = -17069

$ws.Range("H137").Value = 1369.75
$ws.Range("I137").Value = 1170.4546
$ws.Range("J137").Value = 2100.5
$ws.Range("K137").Value = 3511.3638
$ws.Range("L137").Value = 6301.5
$ws.Range("M137").Value = -961.3638000000001
$ws.Range("N137").Value = -11401.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6065.9834
$ws.Range("I32").Value = 5218.125
$ws.Range("K32").Value = 5218.125
$ws.Range("M32").Value = -4931.125

$ws.Range("H57").Value = 12000
$ws.Range("I57").Value = 12000
$ws.Range("K57").Value = 12000
$ws.Range("M57").Value = -11516

$ws.Range("H61").Value = 4660.9375
$ws.Range("I61").Value = 4794.839
$ws.Range("J61").Value = 510
$ws.Range("K61").Value = 4794.839
$ws.Range("L61").Value = 510
$ws.Range("M61").Value = -4582.839
$ws.Range("N61").Value = -934  # new cell

$ws.Range("H122").Value = 2139047.5
$ws.Range("I122").Value = 2139047.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6417142.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6414692.5
$ws.Range("N122").ClearContents()  # was -25180

$ws.Range("H135").Value = 140000
$ws.Range("J135").Value = 140000
$ws.Range("L135").Value = 140000
$ws.Range("N135").Value = -150140

$ws.Range("H136").Value = 4660.9375
$ws.Range("I136").Value = 4794.839
$ws.Range("J136").Value = 510
$ws.Range("K136").Value = 14384.517
$ws.Range("L136").Value = 1530
$ws.Range("M136").Value = -11834.517
$ws.Range("N136").Value = -6630  # new cell

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 40162.75
$ws.Range("J58").Value = 40162.75
$ws.Range("L58").Value = 40162.75
$ws.Range("N58").Value = -40750.75

$ws.Range("H86").Value = 17545722
$ws.Range("I86").Value = 23811208
$ws.Range("K86").Value = 23811208
$ws.Range("M86").Value = -23810085

$ws.Range("H89").Value = 17545722
$ws.Range("I89").Value = 23811208
$ws.Range("K89").Value = 119056040
$ws.Range("M89").Value = -119050424

$ws.Range("H105").Value = 25706
$ws.Range("I105").Value = 101005
$ws.Range("K105").Value = 101005
$ws.Range("M105").Value = -99258

$ws.Range("H134").Value = 5245.647
$ws.Range("I134").Value = 7482.5264
$ws.Range("J134").Value = 2412.2666
$ws.Range("K134").Value = 22447.5792
$ws.Range("L134").Value = 7236.7998
$ws.Range("M134").Value = -19912.5792
$ws.Range("N134").Value = -12306.7998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3291.32
$ws.Range("I31").Value = 1979.1538
$ws.Range("J31").Value = 4712.8335
$ws.Range("K31").Value = 1979.1538
$ws.Range("L31").Value = 4712.8335
$ws.Range("M31").Value = -1684.1538
$ws.Range("N31").Value = -5302.8335

$ws.Range("H34").Value = 3291.32
$ws.Range("I34").Value = 1979.1538
$ws.Range("J34").Value = 4712.8335
$ws.Range("K34").Value = 1979.1538
$ws.Range("L34").Value = 4712.8335
$ws.Range("M34").Value = -1777.1538
$ws.Range("N34").Value = -5116.8335

$ws.Range("H132").Value = 2221.5151
$ws.Range("I132").Value = 1911.1428
$ws.Range("J132").Value = 3959.6
$ws.Range("K132").Value = 5733.428400000001
$ws.Range("L132").Value = 11878.8
$ws.Range("M132").Value = -3203.428400000001
$ws.Range("N132").Value = -16938.8

$ws.Range("H134").Value = 3315.375
$ws.Range("I134").Value = 4102.5625
$ws.Range("J134").Value = 1741
$ws.Range("K134").Value = 12307.6875
$ws.Range("L134").Value = 5223
$ws.Range("M134").Value = -9772.6875
$ws.Range("N134").Value = -10293

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 250472.2
$ws.Range("I5").Value = 611.63635
$ws.Range("J5").Value = 461892.7
$ws.Range("K5").Value = 1834.90905
$ws.Range("L5").Value = 1385678.1
$ws.Range("M5").Value = -1722.90905
$ws.Range("N5").Value = -1385902.1

$ws.Range("H122").Value = 977.1111
$ws.Range("J122").Value = 1157.4
$ws.Range("L122").Value = 10416.6
$ws.Range("N122").Value = -15316.6

$ws.Range("H132").Value = 3448.276
$ws.Range("I132").Value = 3000
$ws.Range("J132").Value = 3456.1404
$ws.Range("K132").Value = 27000
$ws.Range("L132").Value = 31105.2636
$ws.Range("M132").Value = -24470
$ws.Range("N132").Value = -36165.26360000001

$ws.Range("H133").Value = 53683.227
$ws.Range("I133").Value = 172671.83
$ws.Range("J133").Value = 9062.5
$ws.Range("K133").Value = 518015.49
$ws.Range("L133").Value = 27187.5
$ws.Range("M133").Value = -512955.49
$ws.Range("N133").Value = -37307.5

$ws.Range("H135").Value = 250472.2
$ws.Range("I135").Value = 611.63635
$ws.Range("J135").Value = 461892.7
$ws.Range("K135").Value = 5504.72715
$ws.Range("L135").Value = 4157034.3
$ws.Range("M135").Value = -2969.72715
$ws.Range("N135").Value = -4162104.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 9905.714
$ws.Range("I46").Value = 7085.25
$ws.Range("J46").Value = 13666.333
$ws.Range("K46").Value = 7085.25
$ws.Range("L46").Value = 13666.333
$ws.Range("M46").Value = -6929.25
$ws.Range("N46").Value = -13978.333

$ws.Range("H113").Value = 166667600
$ws.Range("I113").Value = 250000690
$ws.Range("J113").Value = 1450
$ws.Range("K113").Value = 250000690
$ws.Range("L113").Value = 1450
$ws.Range("M113").Value = -249998520
$ws.Range("N113").Value = -5790

$ws.Range("H132").Value = 5237.353
$ws.Range("I132").Value = 9087.166999999999
$ws.Range("J132").Value = 3137.4546
$ws.Range("K132").Value = 27261.501
$ws.Range("L132").Value = 9412.363799999999
$ws.Range("M132").Value = -24731.501
$ws.Range("N132").Value = -14472.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 250002930
$ws.Range("I40").Value = 333335900
$ws.Range("K40").Value = 333335900
$ws.Range("M40").Value = -333335764

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2096.25
$ws.Range("I81").Value = 1324.2858
$ws.Range("J81").Value = 7500
$ws.Range("K81").Value = 2648.5716
$ws.Range("L81").Value = 15000
$ws.Range("M81").Value = -1587.5716
$ws.Range("N81").Value = -17122  # new cell

$ws.Range("H84").Value = 2096.25
$ws.Range("I84").Value = 1324.2858
$ws.Range("J84").Value = 7500
$ws.Range("K84").Value = 13242.858
$ws.Range("L84").Value = 75000
$ws.Range("M84").Value = -7938.858
$ws.Range("N84").Value = -85608  # new cell

$ws.Range("H122").Value = 1555.64
$ws.Range("I122").Value = 1354.1818
$ws.Range("K122").Value = 4062.5454
$ws.Range("M122").Value = -1612.5454
